# Replaces the arithmetic answers in the single 20-row x 5-column table.
# Row 11 (1-indexed) also changes which two numbers occupy its first two
# cells (the old 3rd/4th cells are dropped), but the cell COUNT per row
# stays 5 throughout, so every update below is a simple cell-text set
# keyed by (row, column) -- this sidesteps ambiguity from duplicate
# expression text (e.g. "30+26=56" appears twice in the source table)
# and avoids any row/column insert-delete gymnastics.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$row = $t.Rows.Item(1)
$row.Cells.Item(1).Range.Text = "96-61=35"
$row.Cells.Item(2).Range.Text = "78-19=59"
$row.Cells.Item(3).Range.Text = "90-86=4"
$row.Cells.Item(4).Range.Text = "56-37=19"
$row.Cells.Item(5).Range.Text = "47-6=41"

$row = $t.Rows.Item(2)
$row.Cells.Item(1).Range.Text = "56-4=52"
$row.Cells.Item(2).Range.Text = "1+60=61"
$row.Cells.Item(3).Range.Text = "6+89=95"
$row.Cells.Item(4).Range.Text = "93-48=45"
$row.Cells.Item(5).Range.Text = "43-15=28"

$row = $t.Rows.Item(3)
$row.Cells.Item(1).Range.Text = "16+28=44"
$row.Cells.Item(2).Range.Text = "34+55=89"
$row.Cells.Item(3).Range.Text = "43+28=71"
$row.Cells.Item(4).Range.Text = "93+4=97"
$row.Cells.Item(5).Range.Text = "90-64=26"

$row = $t.Rows.Item(4)
$row.Cells.Item(1).Range.Text = "26+67=93"
$row.Cells.Item(2).Range.Text = "57+20=77"
$row.Cells.Item(3).Range.Text = "57+39=96"
$row.Cells.Item(4).Range.Text = "4+53=57"
$row.Cells.Item(5).Range.Text = "57+18=75"

$row = $t.Rows.Item(5)
$row.Cells.Item(1).Range.Text = "18+41=59"
$row.Cells.Item(2).Range.Text = "47-12=35"
$row.Cells.Item(3).Range.Text = "4+59=63"
$row.Cells.Item(4).Range.Text = "34-14=20"
$row.Cells.Item(5).Range.Text = "95-87=8"

$row = $t.Rows.Item(6)
$row.Cells.Item(1).Range.Text = "70+11=81"
$row.Cells.Item(2).Range.Text = "35-17=18"
$row.Cells.Item(3).Range.Text = "44+19=63"
$row.Cells.Item(4).Range.Text = "52+18=70"
$row.Cells.Item(5).Range.Text = "93-44=49"

$row = $t.Rows.Item(7)
$row.Cells.Item(1).Range.Text = "67-21=46"
$row.Cells.Item(2).Range.Text = "82-5=77"
$row.Cells.Item(3).Range.Text = "14+63=77"
$row.Cells.Item(4).Range.Text = "98-65=33"
$row.Cells.Item(5).Range.Text = "97-61=36"

$row = $t.Rows.Item(8)
$row.Cells.Item(1).Range.Text = "35+11=46"
$row.Cells.Item(2).Range.Text = "22+43=65"
$row.Cells.Item(3).Range.Text = "30+56=86"
$row.Cells.Item(4).Range.Text = "76-33=43"
$row.Cells.Item(5).Range.Text = "60-48=12"

$row = $t.Rows.Item(9)
$row.Cells.Item(1).Range.Text = "95-15=80"
$row.Cells.Item(2).Range.Text = "12+2=14"
$row.Cells.Item(3).Range.Text = "56+21=77"
$row.Cells.Item(4).Range.Text = "4+14=18"
$row.Cells.Item(5).Range.Text = "41-0=41"

$row = $t.Rows.Item(10)
$row.Cells.Item(1).Range.Text = "53-9=44"
$row.Cells.Item(2).Range.Text = "88-82=6"
$row.Cells.Item(3).Range.Text = "98-61=37"
$row.Cells.Item(4).Range.Text = "68-49=19"
$row.Cells.Item(5).Range.Text = "95-74=21"

$row = $t.Rows.Item(11)
$row.Cells.Item(1).Range.Text = "40+40=80"
$row.Cells.Item(2).Range.Text = "42+45=87"
$row.Cells.Item(3).Range.Text = "43-4=39"
$row.Cells.Item(4).Range.Text = "25+4=29"
$row.Cells.Item(5).Range.Text = "99-96=3"

$row = $t.Rows.Item(12)
$row.Cells.Item(1).Range.Text = "30+27=57"
$row.Cells.Item(2).Range.Text = "77-27=50"
$row.Cells.Item(3).Range.Text = "40+55=95"
$row.Cells.Item(4).Range.Text = "62-1=61"
$row.Cells.Item(5).Range.Text = "53-20=33"

$row = $t.Rows.Item(13)
$row.Cells.Item(1).Range.Text = "1+77=78"
$row.Cells.Item(2).Range.Text = "21+68=89"
$row.Cells.Item(3).Range.Text = "42-31=11"
$row.Cells.Item(4).Range.Text = "27-10=17"
$row.Cells.Item(5).Range.Text = "36+4=40"

$row = $t.Rows.Item(14)
$row.Cells.Item(1).Range.Text = "91+7=98"
$row.Cells.Item(2).Range.Text = "28-1=27"
$row.Cells.Item(3).Range.Text = "86+1=87"
$row.Cells.Item(4).Range.Text = "14-12=2"
$row.Cells.Item(5).Range.Text = "42+1=43"

$row = $t.Rows.Item(15)
$row.Cells.Item(1).Range.Text = "20+77=97"
$row.Cells.Item(2).Range.Text = "25+63=88"
$row.Cells.Item(3).Range.Text = "3+95=98"
$row.Cells.Item(4).Range.Text = "28+41=69"
$row.Cells.Item(5).Range.Text = "20+48=68"

$row = $t.Rows.Item(16)
$row.Cells.Item(1).Range.Text = "44+52=96"
$row.Cells.Item(2).Range.Text = "30+6=36"
$row.Cells.Item(3).Range.Text = "92-8=84"
$row.Cells.Item(4).Range.Text = "6+52=58"
$row.Cells.Item(5).Range.Text = "22-14=8"

$row = $t.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "1+93=94"
$row.Cells.Item(2).Range.Text = "15-6=9"
$row.Cells.Item(3).Range.Text = "31+24=55"
$row.Cells.Item(4).Range.Text = "89-62=27"
$row.Cells.Item(5).Range.Text = "60+27=87"

$row = $t.Rows.Item(18)
$row.Cells.Item(1).Range.Text = "14+84=98"
$row.Cells.Item(2).Range.Text = "12+42=54"
$row.Cells.Item(3).Range.Text = "21+44=65"
$row.Cells.Item(4).Range.Text = "25-21=4"
$row.Cells.Item(5).Range.Text = "2+94=96"

$row = $t.Rows.Item(19)
$row.Cells.Item(1).Range.Text = "51-21=30"
$row.Cells.Item(2).Range.Text = "2+71=73"
$row.Cells.Item(3).Range.Text = "95-61=34"
$row.Cells.Item(4).Range.Text = "62-9=53"
$row.Cells.Item(5).Range.Text = "77-7=70"

$row = $t.Rows.Item(20)
$row.Cells.Item(1).Range.Text = "25+20=45"
$row.Cells.Item(2).Range.Text = "85-48=37"
$row.Cells.Item(3).Range.Text = "37+60=97"
$row.Cells.Item(4).Range.Text = "99-93=6"
$row.Cells.Item(5).Range.Text = "0+39=39"
